$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all touched cells so numeric-looking strings
# (e.g. "592.08", "0.530", "67.842.69") are preserved exactly as text
# rather than being auto-converted to numbers by Excel.

$ws.Range('D2:E2').NumberFormat = '@'
$ws.Range('D2').Value = '67.842.69'
$ws.Range('E2').Value = '  +2.38%  '

$ws.Range('D3:E3').NumberFormat = '@'
$ws.Range('D3').Value = '2.528.39'
$ws.Range('E3').Value = '  -1.10%  '

$ws.Range('D4:E4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5:E5').NumberFormat = '@'
$ws.Range('D5').Value = '592.08'
$ws.Range('E5').Value = '  +1.76%  '

$ws.Range('D6:E6').NumberFormat = '@'
$ws.Range('D6').Value = '176.37'
$ws.Range('E6').Value = '  +5.76%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('D8:E8').NumberFormat = '@'
$ws.Range('D8').Value = '0.530'
$ws.Range('E8').Value = '  +0.78%  '

$ws.Range('D9:E9').NumberFormat = '@'
$ws.Range('D9').Value = '2.525.85'
$ws.Range('E9').Value = '  -1.16%  '

$ws.Range('D10:E10').NumberFormat = '@'
$ws.Range('D10').Value = '0.142'
$ws.Range('E10').Value = '  +1.79%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.08%  '

$ws.Range('D12:E12').NumberFormat = '@'
$ws.Range('D12').Value = '5.16'
$ws.Range('E12').Value = '  +0.29%  '

$ws.Range('D13:E13').NumberFormat = '@'
$ws.Range('D13').Value = '0.344'
$ws.Range('E13').Value = '  -2.53%  '

$ws.Range('D14:E14').NumberFormat = '@'
$ws.Range('D14').Value = '26.82'
$ws.Range('E14').Value = '  +0.83%  '

$ws.Range('D15:E15').NumberFormat = '@'
$ws.Range('D15').Value = '2.985.41'
$ws.Range('E15').Value = '  -1.40%  '

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.75%  '

$ws.Range('D17:E17').NumberFormat = '@'
$ws.Range('D17').Value = '67.606.42'
$ws.Range('E17').Value = '  +2.02%  '

$ws.Range('D18:E18').NumberFormat = '@'
$ws.Range('D18').Value = '2.496.41'
$ws.Range('E18').Value = '  -2.66%  '

$ws.Range('D19:E19').NumberFormat = '@'
$ws.Range('D19').Value = '8.05'
$ws.Range('E19').Value = '  +4.51%  '

$ws.Range('D20:E20').NumberFormat = '@'
$ws.Range('D20').Value = '11.49'
$ws.Range('E20').Value = '  +1.08%  '

$ws.Range('D21:E21').NumberFormat = '@'
$ws.Range('D21').Value = '359.18'
$ws.Range('E21').Value = '  +3.05%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.31%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.54%  '

$ws.Range('D24:E24').NumberFormat = '@'
$ws.Range('D24').Value = '1.99'
$ws.Range('E24').Value = '  +4.22%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.00%  '

$ws.Range('B26:E26').NumberFormat = '@'
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '70.93'
$ws.Range('E26').Value = '  +2.55%  '

$ws.Range('B27:E27').NumberFormat = '@'
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = '10.29'
$ws.Range('E27').Value = '  +3.89%  '

$ws.Range('B28:E28').NumberFormat = '@'
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '0.996'
$ws.Range('E28').Value = '  -0.40%  '

$ws.Range('B29:E29').NumberFormat = '@'
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '2.654.25'
$ws.Range('E29').Value = '  -1.31%  '

$ws.Range('D30:E30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0989'
$ws.Range('E30').Value = '  +0.36%  '

$ws.Range('D31:E31').NumberFormat = '@'
$ws.Range('D31').Value = '558.94'
$ws.Range('E31').Value = '  +6.39%  '

$ws.Range('D32:E32').NumberFormat = '@'
$ws.Range('D32').Value = '8.30'
$ws.Range('E32').Value = '  +0.50%  '

$ws.Range('D33:E33').NumberFormat = '@'
$ws.Range('D33').Value = '1.36'
$ws.Range('E33').Value = '  +3.02%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.02%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.16%  '

$ws.Range('D36:E36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.20%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.97%  '

$ws.Range('D38:E38').NumberFormat = '@'
$ws.Range('D38').Value = '156.20'
$ws.Range('E38').Value = '  -0.69%  '

$ws.Range('D39:E39').NumberFormat = '@'
$ws.Range('D39').Value = '18.76'
$ws.Range('E39').Value = '  +0.37%  '

$ws.Range('D40:E40').NumberFormat = '@'
$ws.Range('D40').Value = '18.60'
$ws.Range('E40').Value = '  +1.60%  '

$ws.Range('D41:E41').NumberFormat = '@'
$ws.Range('D41').Value = '1.82'
$ws.Range('E41').Value = '  +3.00%  '

$ws.Range('D42:E42').NumberFormat = '@'
$ws.Range('D42').Value = '0.356'
$ws.Range('E42').Value = '  -0.14%  '

$ws.Range('D43:E43').NumberFormat = '@'
$ws.Range('D43').Value = '5.17'
$ws.Range('E43').Value = '  +1.98%  '

$ws.Range('D44:E44').NumberFormat = '@'
$ws.Range('D44').Value = '2.53'
$ws.Range('E44').Value = '  +4.90%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.03%  '

$ws.Range('D46:E46').NumberFormat = '@'
$ws.Range('D46').Value = '147.68'
$ws.Range('E46').Value = '  -0.23%  '

$ws.Range('D47:E47').NumberFormat = '@'
$ws.Range('D47').Value = '0.560'
$ws.Range('E47').Value = '  -0.27%  '

$ws.Range('D48:E48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0278'
$ws.Range('E48').Value = '  -2.54%  '

$ws.Range('D49:E49').NumberFormat = '@'
$ws.Range('D49').Value = '3.71'
$ws.Range('E49').Value = '  +0.24%  '

$ws.Range('D50:E50').NumberFormat = '@'
$ws.Range('D50').Value = '1.69'
$ws.Range('E50').Value = '  -1.63%  '

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.42%  '
